$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Jan -> Weasel, Unknown Sex Count (J) = 1
$ws.Range("F3").Value = "Weasel"
$ws.Range("J3").Value = 1

# Row 4: Weasel -> Marten, Female Count (I) = 1, Unknown Sex Count (J) = 0
$ws.Range("F4").Value = "Marten"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0

# Row 5: Marten -> Mink, Female Count (I) = 1
$ws.Range("F5").Value = "Mink"
$ws.Range("I5").Value = 1

# Row 6: January -> February, Mink -> Lynx, Female Count (I) = 1
$ws.Range("E6").Value = "February"
$ws.Range("F6").Value = "Lynx"
$ws.Range("I6").Value = 1

# Remove the now-obsolete trailing rows (old row 7 blank-February row and
# old row 8 Lynx row); their data has been folded into row 6 above.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
